$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Fecha 44280 -> 44315, Calidad Primera -> Especial, Volumen 30 -> 10,
#        precios 25000 -> 30000, unidad "18 kilos empedrada" -> "20 kilos empedrada",
#        precio/kg 1389 -> 1500, kilos 18 -> 20
$ws.Range("D5").Value = 44315
$ws.Range("I5").Value = "Especial"
$ws.Range("J5").Value = 10
$ws.Range("K5").Value = 30000
$ws.Range("L5").Value = 30000
$ws.Range("M5").Value = 30000
$ws.Range("N5").Value = "`$/caja 20 kilos empedrada"
$ws.Range("P5").Value = 1500
$ws.Range("Q5").Value = 20

# Row 6: Fecha 44285 -> 44315, precios 25000 -> 15000,
#        unidad "18 kilos empedrada" -> "15 kilos granel",
#        precio/kg 1389 -> 1000, kilos 18 -> 15
$ws.Range("D6").Value = 44315
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 15000
$ws.Range("N6").Value = "`$/caja 15 kilos granel"
$ws.Range("P6").Value = 1000
$ws.Range("Q6").Value = 15

# Row 7: Fecha 44315 -> 44280, Calidad Especial -> Primera, Volumen 10 -> 30,
#        precios 30000 -> 25000, unidad "20 kilos empedrada" -> "18 kilos empedrada",
#        precio/kg 1500 -> 1389, kilos 20 -> 18
$ws.Range("D7").Value = 44280
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 30
$ws.Range("K7").Value = 25000
$ws.Range("L7").Value = 25000
$ws.Range("M7").Value = 25000
$ws.Range("N7").Value = "`$/caja 18 kilos empedrada"
$ws.Range("P7").Value = 1389
$ws.Range("Q7").Value = 18

# Row 8: Fecha 44315 -> 44285, precios 15000 -> 25000,
#        unidad "15 kilos granel" -> "18 kilos empedrada",
#        precio/kg 1000 -> 1389, kilos 15 -> 18
$ws.Range("D8").Value = 44285
$ws.Range("K8").Value = 25000
$ws.Range("L8").Value = 25000
$ws.Range("M8").Value = 25000
$ws.Range("N8").Value = "`$/caja 18 kilos empedrada"
$ws.Range("P8").Value = 1389
$ws.Range("Q8").Value = 18
